$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row: drop first_name/last_name, move ip_address to column B,
# add server_port and online_status columns.
$ws.Range("A1").Value = "user_id"
$ws.Range("B1").Value = "ip_address"
$ws.Range("C1").Value = "server_port"
$ws.Range("D1").Value = "online_status"

# Capture the existing ip_address values (currently in column D) before
# column B (first_name) gets overwritten, so we don't clobber data we
# still need to read.
$ipAddresses = @()
for ($r = 2; $r -le 13; $r++) {
    $ipAddresses += $ws.Cells.Item($r, 4).Value()
}

# Rewrite each data row: keep user_id in A, move the ip address into B,
# and populate the new server_port / online_status numeric columns.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 2).Value = $ipAddresses[$r - 2]
    $ws.Cells.Item($r, 3).Value = 1200
    $ws.Cells.Item($r, 4).Value = 0
}

# Update the active selection to D2, matching the saved view state.
$ws.Range("D2").Select()
